$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.543.93"
$ws.Range("E2").Value = "  +5.39%  "

$ws.Range("D3").Value = "2.297.91"
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.76"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.25"
$ws.Range("E6").Value = "  +11.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  +1.78%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  +5.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.53"
$ws.Range("E10").Value = "  +9.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.46"
$ws.Range("E12").Value = "  +6.73%  "

$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").Value = "2.646.18"
$ws.Range("E14").Value = "  +3.15%  "

$ws.Range("D15").Value = "2.295.22"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.83"
$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("E17").Value = "  +4.44%  "

$ws.Range("D18").Value = "46.521.75"
$ws.Range("E18").Value = "  +5.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  +6.26%  "

$ws.Range("E20").Value = "  +3.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.17"
$ws.Range("E22").Value = "  +3.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.79"
$ws.Range("E23").Value = "  +5.41%  "

$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("E26").Value = "  +4.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.40"
$ws.Range("E27").Value = "  +8.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.91"
$ws.Range("E29").Value = "  +5.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.03"
$ws.Range("E30").Value = "  +4.05%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.84"
$ws.Range("E31").Value = "  +13.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.64"
$ws.Range("E32").Value = "  +2.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.08"
$ws.Range("E33").Value = "  -3.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0793"
$ws.Range("E34").Value = "  +3.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.27"
$ws.Range("E35").Value = "  +14.44%  "

$ws.Range("E36").Value = "  +8.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("E37").Value = "  +0.82%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  +4.82%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.01"
$ws.Range("E39").Value = "  +17.82%  "

$ws.Range("E40").Value = "  +10.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.36"
$ws.Range("E41").Value = "  +6.18%  "

$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +10.56%  "

$ws.Range("D45").Value = "1.814.29"
$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.50"
$ws.Range("E46").Value = "  +20.73%  "

$ws.Range("E47").Value = "  +4.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.33"
$ws.Range("E48").Value = "  +7.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.89"
$ws.Range("E49").Value = "  +5.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.92"
$ws.Range("E50").Value = "  +1.27%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.94"
$ws.Range("E51").Value = "  +2.34%  "

